# Add the new week "24/01/2022 - 30/01/2022" data to all three sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Asl Sorveglianza": Data | Asl Sorveglianza | N. Positivi
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Asl Sorveglianza")

$ws1.Range("A107").Value = "24/01/2022 - 30/01/2022"
$ws1.Range("B107").Value = "AZIENDA USL TOSCANA SUD-EST"
$ws1.Range("C107").Value = 33

$ws1.Range("A108").Value = "24/01/2022 - 30/01/2022"
$ws1.Range("B108").Value = "AZIENDA USL TOSCANA CENTRO"
$ws1.Range("C108").Value = 152

$ws1.Range("A109").Value = "24/01/2022 - 30/01/2022"
$ws1.Range("B109").Value = "AZIENDA USL TOSCANA NORD-OVEST"
$ws1.Range("B109").Font.Color = 0
$ws1.Range("C109").Value = 94

$ws1.Range("B110").Value = "Totale"

# ---------------------------------------------------------------------------
# Sheet "Professione": Data | Professione | N. Positivi | Totale contatti
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Professione")

$ws2.Range("A78").Value = "24/01/2022 - 30/01/2022"
$ws2.Range("B78").Value = "Insegnante"
$ws2.Range("C78").Value = 264
$ws2.Range("D78").Value = 2561
$ws2.Range("D78").NumberFormat = "#,##0"

$ws2.Range("A79").Value = "24/01/2022 - 30/01/2022"
$ws2.Range("B79").Value = "Personale non docente"
$ws2.Range("C79").Value = 15
$ws2.Range("D79").Value = 112
$ws2.Range("D79").NumberFormat = "#,##0"

$ws2.Range("B80").Value = "Totale"
$ws2.Range("C80").Value = 279
$ws2.Range("D80").Value = 2673
$ws2.Range("D80").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Sheet "Sesso ed età": Data | Età | Sesso | Totale
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sesso ed età")

$ws3.Range("A139").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B139").Value = "0-18"
$ws3.Range("C139").Value = "F"
$ws3.Range("D139").Value = 2

$ws3.Range("A140").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B140").Value = "0-18"
$ws3.Range("C140").Value = "M"
$ws3.Range("D140").Value = 1

$ws3.Range("A141").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B141").Value = "19-34"
$ws3.Range("C141").Value = "F"
$ws3.Range("D141").Value = 32

$ws3.Range("A142").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B142").Value = "19-34"
$ws3.Range("B142").Font.Color = 0
$ws3.Range("C142").Value = "M"
$ws3.Range("D142").Value = 7

$ws3.Range("A143").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B143").Value = "35-49"
$ws3.Range("C143").Value = "F"
$ws3.Range("D143").Value = 135

$ws3.Range("A144").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B144").Value = "35-49"
$ws3.Range("C144").Value = "M"
$ws3.Range("D144").Value = 16

$ws3.Range("A145").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B145").Value = "50-64"
$ws3.Range("C145").Value = "F"
$ws3.Range("D145").Value = 77

$ws3.Range("A146").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B146").Value = "50-64"
$ws3.Range("C146").Value = "M"
$ws3.Range("D146").Value = 7

$ws3.Range("A147").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B147").Value = "65-79"
$ws3.Range("C147").Value = "F"
$ws3.Range("D147").Value = 1

$ws3.Range("A148").Value = "24/01/2022 - 30/01/2022"
$ws3.Range("B148").Value = "65-79"
$ws3.Range("C148").Value = "M"
$ws3.Range("D148").Value = 1

# ---------------------------------------------------------------------------
# Restore selection / active sheet to match the final saved view state.
# ---------------------------------------------------------------------------
$ws1.Select()
$ws1.Range("C110").Select()

$ws2.Select()
$ws2.Range("A78").Select()

$ws3.Select()
$ws3.Range("D149").Select()
